$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: user overwrites row 2 (was "Exotic Scimitar") while drafting a new weapon "sword_33"
$ws.Range("B2").Value = "sword_33"
$ws.Range("C2").Value = "sword_33 (1)"
$ws.Range("G2").Value = "Onyx"
$ws.Range("H2").Value = "Iron"

# Step 2: user adds a new row 3 for weapon "sword_32"
$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("B3").Value = "sword_32"
$ws.Range("C3").Value = "sword_32 (1)"
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = "Onyx"
$ws.Range("H3").Value = "Bloodstone"

# Step 3: user re-adds the original "Exotic Scimitar" entry as new row 4
$ws.Range("A2").Copy($ws.Range("A4"))
$ws.Range("B4").Value = "Exotic Scimitar"
$ws.Range("C4").Value = "Exotic Scimitar (1)"
$ws.Range("D4").Value = 11
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 22
$ws.Range("G4").Value = "Silver"
$ws.Range("H4").Value = "Zircon"
